$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 408, pushing existing rows 408-485 down to 412-489.
$ws.Rows("408:411").Insert()

# New row data (Sandia, Vega Central Mapocho de Santiago, Metropolitana region, date 44637 = 2022-03-17)
$newRows = @(
    @{ Row = 408; I = "Extra";   J = 160; K = 3200; L = 3500; M = 3350; O = "Región Metropolitana"; P = 3350 },
    @{ Row = 409; I = "Primera"; J = 340; K = 2800; L = 3000; M = 2900; O = "Región Metropolitana"; P = 2900 },
    @{ Row = 410; I = "Segunda"; J = 250; K = 2400; L = 2600; M = 2500; O = "Región Metropolitana"; P = 2500 },
    @{ Row = 411; I = "Tercera"; J = 61;  K = 2000; L = 2000; M = 2000; O = "Región Metropolitana"; P = 2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44637
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112028
    $ws.Cells.Item($row, 7).Value = "Sandia"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
